$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.223.24'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.773.04'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.52%  '
$ws.Range('E3').ClearFormats()
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.10'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.54%  '
$ws.Range('E5').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5260'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +11.26%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3678'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +6.82%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.71'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E9').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('E10').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.01%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.39'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.74%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.049'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.34%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.767.89'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.15%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.917'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.62'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001044'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06421'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('E19').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E20').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.798'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.64%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.278.94'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E23').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.33%  '
$ws.Range('E24').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.87'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.11'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('E27').ClearFormats()
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B28').ClearFormats()
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C28').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.971.00'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('E28').ClearFormats()
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('B29').ClearFormats()
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C29').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.319'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +11.07%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.63'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.058'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.44%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09764'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.96%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.554'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.81%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.617'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02235'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05960'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('E36').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E37').ClearFormats()
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('B38').ClearFormats()
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C38').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.825'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.77%  '
$ws.Range('E38').ClearFormats()
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('B39').ClearFormats()
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C39').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6126'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.39%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2017'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.428'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.055'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.80%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.134'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.03%  '
$ws.Range('E43').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5743'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('E45').ClearFormats()
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('B46').ClearFormats()
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C46').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.622'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.17'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.54%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.873'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.112'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.37%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06692'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.91%  '
$ws.Range('E50').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.08%  '
$ws.Range('E51').ClearFormats()
